$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-98 down to 58-99.
$ws.Rows("57:57").Insert()

# Populate the newly inserted row 57 with the new record's data.
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = "Vega Monumental Concepción"
$ws.Range("C57").Value = "Bíobío"
$ws.Range("D57").Value = 44673
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 100112021
$ws.Range("G57").Value = "Ají"
$ws.Range("H57").Value = "Inferno"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 22
$ws.Range("K57").Value = 16000
$ws.Range("L57").Value = 17000
$ws.Range("M57").Value = 16545
$ws.Range("N57").Value = "$/caja 12 kilos"
$ws.Range("O57").Value = "Región de Arica y Parinacota"
$ws.Range("P57").Value = 1379
$ws.Range("Q57").Value = 12
$ws.Range("R57").Value = "Hortaliza"
